$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename/reorder column headers ---
$ws.Range("A1").Value = "YEAR"
$ws.Range("B1").Value = "tuitionandfees"
$ws.Range("C1").Value = "booksAndSupplies"
$ws.Range("D1").Value = "roomAndBoard"
$ws.Range("E1").Value = "other"
$ws.Range("F1").Value = "total"

# --- New "total" column: F2 standalone formula, F3:F14 filled/shared ---
$ws.Range("F2").Formula = "=SUM(B2:E2)"
$ws.Range("F3:F14").Formula = "=SUM(B3:E3)"

# --- Update selection to match the saved view ---
[void]$ws.Range("G20").Select()
